$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("G1").Value = "Precio Oro"
$ws.Range("H1").Value = "Precio Hierro"
$ws.Range("I1").Value = "Precio Comida"
$ws.Range("J1").Value = "Consumo Comida"

$ws.Range("G1:J1").EntireColumn.AutoFit() | Out-Null

$ws.Range("H3").Select()
